# Update the dashboard export so the data reflects B450 (instead of U200),
# and correct the Plan values (column D) for several days per the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label for the Plan/Actual/Status(+/-) group (merged D4:F4): U200 -> B450
$ws.Range("D4").Value = "B450"

# Corrected "Plan" values in column D
$ws.Range("D7").Value = 150
$ws.Range("D11").Value = 115
$ws.Range("D12").Value = 150
$ws.Range("D13").Value = 115
$ws.Range("D14").Value = 150
$ws.Range("D18").Value = 115
$ws.Range("D19").Value = 150
$ws.Range("D20").Value = 115
$ws.Range("D21").Value = 150
$ws.Range("D25").Value = 115
$ws.Range("D26").Value = 150
$ws.Range("D27").Value = 115
$ws.Range("D28").Value = 150
$ws.Range("D32").Value = 115
$ws.Range("D34").Value = 0
